$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3-5 (delete from bottom up to keep row indices stable)
$ws.Range("A5:B5").EntireRow.Delete()
$ws.Range("A4:B4").EntireRow.Delete()
$ws.Range("A3:B3").EntireRow.Delete()

# Update remaining row 2 values (date + temperature)
$ws.Range("A2").Value = 44835.125
$ws.Range("B2").Value = 12.4
